# Add the 2021 season rows (114-125) to "Historical Stats - Money League"
# Sheet1, extending the existing data table and its shared "champion" flag
# formula in column J, then move the active selection to H117 to match
# where the author was working after the paste.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each tuple: Year, Manager, Finish, DraftPosition, Points, PointsAllowed, Playoffs, Wins, Bye
$newRows = @(
    @(114, 2021, "Colin",    1,  8, 1874.4,  1833.72, 1, 9,  0),
    @(115, 2021, "John",     2, 10, 2131.22, 1822.68, 1, 7,  1),
    @(116, 2021, "Charles",  3,  8, 2029.34, 1887.08, 1, 12, 0),
    @(117, 2021, "Jennifer", 4, 11, 2230.96, 1814.1,  1, 8,  1),
    @(118, 2021, "Chester",  5, 10, 1980.48, 1788.14, 1, 1,  0),
    @(119, 2021, "EricR",    6,  9, 1991.1,  1931.18, 1, 10, 0),
    @(120, 2021, "Mike",     7,  5, 1656.64, 1876.48, 0, 5,  0),
    @(121, 2021, "EricNC",   8,  7, 1897.58, 1875.08, 0, 11, 0),
    @(122, 2021, "ChrisNC",  9,  3, 1703,    1982.58, 0, 4,  0),
    @(123, 2021, "Erik",    10,  6, 1634.9,  1847.78, 0, 3,  0),
    @(124, 2021, "Marcus",  11,  3, 1655.42, 1955.62, 0, 6,  0),
    @(125, 2021, "Alex",    12,  4, 1720.02, 1890.62, 0, 2,  0)
)

foreach ($row in $newRows) {
    $r          = $row[0]
    $year       = $row[1]
    $manager    = $row[2]
    $finish     = $row[3]
    $draftPos   = $row[4]
    $points     = $row[5]
    $ptsAllowed = $row[6]
    $playoffs   = $row[7]
    $wins       = $row[8]
    $bye        = $row[9]

    $ws.Cells.Item($r, 1).Value = $year
    $ws.Cells.Item($r, 2).Value = $manager
    $ws.Cells.Item($r, 3).Value = $finish
    $ws.Cells.Item($r, 4).Value = $draftPos
    $ws.Cells.Item($r, 5).Value = $points
    $ws.Cells.Item($r, 6).Value = $ptsAllowed
    $ws.Cells.Item($r, 7).Value = $playoffs
    $ws.Cells.Item($r, 8).Value = $wins
    $ws.Cells.Item($r, 9).Value = $bye
    $ws.Cells.Item($r, 10).Formula = "=IF(C$r=1,1,0)"

    # Match the surrounding table's cell style (center-aligned), which is
    # what style index "1" represents in this workbook.
    $ws.Range("A" + $r + ":J" + $r).HorizontalAlignment = -4108
}

# Recalculate so the new J-column formulas carry a cached value, and move
# the selection the way it ends up after the paste/scroll in the diff.
$wb.Application.Calculate()
$ws.Range("H117").Select()
